# Move to new sample data: shift the recurring-transaction dates from the
# 6th/13th-of-the-month cadence to the 11th/23rd-of-the-month cadence, and
# insert the new $4000 "AA__TEST__3" transaction (previously row 18) into
# its correct chronological slot (row 13, 2022-06-18), shuffling the
# transactions that used to occupy rows 14-18 down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaction")

# --- Rows 2-12: only the date (column B) moves from 6th/13th -> 11th/23rd ---
$ws.Range("B2").Value  = 44572
$ws.Range("B3").Value  = 44584
$ws.Range("B4").Value  = 44603
$ws.Range("B5").Value  = 44615
$ws.Range("B6").Value  = 44631
$ws.Range("B7").Value  = 44643
$ws.Range("B8").Value  = 44662
$ws.Range("B9").Value  = 44674
$ws.Range("B10").Value = 44692
$ws.Range("B11").Value = 44704
$ws.Range("B12").Value = 44723

# --- Row 13: becomes the $4000 "AA__TEST__3" transaction (2022-06-18) ---
$ws.Range("B13").Value = 44730
$ws.Range("C13").Value = "AA__TEST__3"
$ws.Range("D13").Value = -4000
$ws.Range("E13").Value = "AA:__TEST__:D"
$ws.Range("O13").Value = "AA__TEST__3"

# --- Row 14: becomes an "AA__TEST__2" transaction (2022-06-23) ---
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = 44735
$ws.Range("C14").Value = "AA__TEST__2"
$ws.Range("E14").Value = "AA:__TEST__:C"
$ws.Range("M14").Value = $false
$ws.Range("O14").Value = "AA__TEST__2"

# --- Row 15: becomes an "AA__TEST__1" transaction (2022-07-11) ---
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = 44753
$ws.Range("C15").Value = "AA__TEST__1"
$ws.Range("D15").Value = -250.00
$ws.Range("E15").ClearContents()
$ws.Range("M15").Value = $true
$ws.Range("O15").Value = "AA__TEST__1"

# --- Row 16: becomes an "AA__TEST__2" transaction (2022-07-23) ---
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 44765
$ws.Range("C16").Value = "AA__TEST__2"
$ws.Range("D16").Value = -250
$ws.Range("E16").Value = "AA:__TEST__:C"
$ws.Range("M16").Value = $false
$ws.Range("O16").Value = "AA__TEST__2"

# --- Row 17: becomes an "AA__TEST__1" transaction (2022-08-11) ---
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = 44784
$ws.Range("C17").Value = "AA__TEST__1"
$ws.Range("D17").Value = -250.00
$ws.Range("E17").ClearContents()
$ws.Range("M17").Value = $true
$ws.Range("O17").Value = "AA__TEST__1"

# --- Row 18: becomes an "AA__TEST__2" transaction (2022-08-23) ---
$ws.Range("B18").Value = 44796
$ws.Range("C18").Value = "AA__TEST__2"
$ws.Range("D18").Value = -250
$ws.Range("E18").Value = "AA:__TEST__:C"
$ws.Range("O18").Value = "AA__TEST__2"

# --- Rows 19-26: only the date (column B) moves from 6th/13th -> 11th/23rd ---
$ws.Range("B19").Value = 44815
$ws.Range("B20").Value = 44827
$ws.Range("B21").Value = 44845
$ws.Range("B22").Value = 44857
$ws.Range("B23").Value = 44876
$ws.Range("B24").Value = 44888
$ws.Range("B25").Value = 44906
$ws.Range("B26").Value = 44918
